# "Planning previsionnel" update
#
# Reworks the task-ownership column ("Groupe", column A) for the project
# plan on sheet "projet":
#   - Splits the previously-merged A18:A23 block into individually
#     labelled rows (Alexandre / Bruno / Marc / Groupe / Groupe), adding
#     two brand-new names ("Bruno", "Marc") that did not exist before.
#   - Splits the previously-merged A24:A26 block into individually
#     labelled rows (all "Alexandre").
#   - Removes the (now redundant, since rows are no longer grouped by a
#     shared top border) top border from the whole A9:A27 owner column,
#     and makes every cell there uniformly center/center aligned.
#   - Updates several PLAN/REEL duration cells (columns E/F) and the
#     resulting %-achieved cells (column G) for a handful of tasks.
#   - Leaves the two leftover helper cells (A28/A29) fully blank again.
#   - Restores the selection to the still-merged A13:A17 block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("projet")

# ---------------------------------------------------------------------
# 1) Duration / completion numbers (columns E, F, G) for rows 11-16 & 24
# ---------------------------------------------------------------------
$ws.Range("G11").Value = 0.8

$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5

$ws.Range("E13").ClearContents()
$ws.Range("F13").ClearContents()
$ws.Range("G13").Value = 0.2

$ws.Range("E14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("G14").Value = 0.9

$ws.Range("E15").ClearContents()
$ws.Range("F15").ClearContents()
$ws.Range("G15").Value = 0

$ws.Range("E16").ClearContents()
$ws.Range("F16").ClearContents()
$ws.Range("G16").Value = 0

$ws.Range("E24").ClearContents()
$ws.Range("F24").ClearContents()
$ws.Range("G24").Value = 0

# ---------------------------------------------------------------------
# 2) Break up the A18:A23 "Groupe" merge into individually-named owners
# ---------------------------------------------------------------------
$ws.Range("A18:A23").UnMerge()
$ws.Range("A18").ClearContents()
$ws.Range("A19").Value = "Alexandre"
$ws.Range("A20").Value = "Bruno"
$ws.Range("A21").Value = "Marc"
$ws.Range("A22").Value = "Groupe"
$ws.Range("A23").Value = "Groupe"

# ---------------------------------------------------------------------
# 3) Break up the A24:A26 "Alexandre" merge the same way
# ---------------------------------------------------------------------
$ws.Range("A24:A26").UnMerge()
$ws.Range("A25").Value = "Alexandre"
$ws.Range("A26").Value = "Alexandre"

# ---------------------------------------------------------------------
# 4) Drop the top border running across the owner column and make every
#    row center-aligned both ways (previously some rows were only
#    vertically centered)
# ---------------------------------------------------------------------
$ownerCol = $ws.Range("A9:A27")
$ownerCol.Borders.Item(12).LineStyle = -4142
$ownerCol.Borders.Item(8).LineStyle = -4142
$ownerCol.HorizontalAlignment = -4108
$ownerCol.VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 5) The two former spacer cells go back to being fully empty/unstyled
# ---------------------------------------------------------------------
$ws.Range("A28").Clear()
$ws.Range("A29").Clear()

# ---------------------------------------------------------------------
# 6) Restore view/selection on the still-merged A13:A17 block
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 8
$ws.Range("A13:A17").Select()
